# Recolor the "Código" (ID), "Historia de Usuario" (title) and
# "Criterios de Aceptación" (description) cells of several user-story
# rows to the same green (92D050) already used by the earlier rows in
# the "Usuarios (Buscadores)" table.
#
# OLE color values are 0x00BBGGRR, so RGB 92D050 -> 0x50D092 -> 5296274.

$d = $word.ActiveDocument
$green = 5296274

# Table 1 = "Usuarios (Buscadores)": header row + USU-01..USU-09.
# USU-01..USU-04 are already green; USU-05..USU-09 are rows 6..10.
$usuarios = $d.Tables.Item(1)
for ($r = 6; $r -le 10; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $usuarios.Cell($r, $c).Range.Font.Color = $green
    }
}

# Table 2 = "Panel de Administración": header row + ADM-01..ADM-05.
# Only ADM-01 and ADM-02 (rows 2..3) change color.
$admin = $d.Tables.Item(2)
for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $admin.Cell($r, $c).Range.Font.Color = $green
    }
}

# Table 4 = "Seguridad y Autenticación": header row + SEG-01..SEG-03.
$seguridad = $d.Tables.Item(4)
for ($r = 2; $r -le 4; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $seguridad.Cell($r, $c).Range.Font.Color = $green
    }
}

# Table 5 = "Rubros (Administrador)": header row + RUB-01..RUB-05.
$rubros = $d.Tables.Item(5)
for ($r = 2; $r -le 6; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $rubros.Cell($r, $c).Range.Font.Color = $green
    }
}
